$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # Spaghetti
$ws3 = $wb.Worksheets.Item(3)   # MER

# ------------------------------------------------------------------
# Sheet "Spaghetti" (sheet2)
# ------------------------------------------------------------------

# Remove the scratch calculations that used to live in J3:M3
$ws2.Range("J3:M3").ClearContents()

# Remove the old "Cost" header in I2 (column no longer used here)
$ws2.Range("I2").ClearContents()

# Correct value used for hex 8 heat duty
$ws2.Range("B11").Value = 51.15

# New "HW1 values" comparison column next to the HEX table
$ws2.Range("I13").Value = "HW1 values"
$ws2.Range("I13").Font.Bold = $true

$ws2.Range("I14").Value = 39.01

$ws2.Range("I15").Value = 259503.7
$ws2.Range("J15").Value = "number of hex = 6"

$ws2.Range("I16").Value = 8160

$ws2.Range("I17").Value = 73036

# New "Cost/Area" summary row
$ws2.Range("G18").Value = "Cost/Area"
$ws2.Range("G18").Font.Bold = $true
$ws2.Range("H18").Formula = "=H17/H14"
$ws2.Range("I18").Formula = "=I17/I14"

$ws2.Range("J3").Select()

# ------------------------------------------------------------------
# Sheet "MER" (sheet3)
# ------------------------------------------------------------------

# New "HW1 values" comparison column next to the HEX table
$ws3.Range("I8").Value = "HW1 values"
$ws3.Range("I8").Font.Bold = $true

$ws3.Range("I9").Value = 39.01

# Insert a new "Cost Hex" row before the old "Utility Costs" row
$ws3.Rows.Item(10).Insert()

$ws3.Range("G10").Value = "Cost Hex"
$ws3.Range("G10").Font.Bold = $true
$ws3.Range("H10").Formula = "=6*(40000+500*H9/6)"
$ws3.Range("I10").Value = 259503.7
$ws3.Range("J10").Value = "number of hex = 6"

$ws3.Range("I11").Value = 8160

$ws3.Range("I12").Value = 73036

# New "Cost/Area" summary row
$ws3.Range("G13").Value = "Cost/Area"
$ws3.Range("G13").Font.Bold = $true
$ws3.Range("H13").Formula = "=H12/H9"
$ws3.Range("I13").Formula = "=I12/I9"

$ws3.Range("H19").Select()
